$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.200.17'
$ws.Range('E2').Value = '  -5.68%  '
$ws.Range('D3').Value = '2.558.50'
$ws.Range('E3').Value = '  -1.29%  '
$ws.Range('D5').Value = "'300.20"
$ws.Range('E5').Value = '  -2.67%  '
$ws.Range('D6').Value = "'92.99"
$ws.Range('E6').Value = '  -5.55%  '
$ws.Range('D7').Value = "'0.577"
$ws.Range('E7').Value = '  -2.78%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  -4.19%  '
$ws.Range('D10').Value = "'36.02"
$ws.Range('E10').Value = '  -6.43%  '
$ws.Range('E11').Value = '  -3.52%  '
$ws.Range('D12').Value = "'7.77"
$ws.Range('E12').Value = '  -3.73%  '
$ws.Range('D13').Value = "'0.109"
$ws.Range('E13').Value = '  +1.89%  '
$ws.Range('D14').Value = '2.947.80'
$ws.Range('E14').Value = '  -1.50%  '
$ws.Range('D15').Value = '2.566.19'
$ws.Range('E15').Value = '  -0.82%  '
$ws.Range('D16').Value = "'0.877"
$ws.Range('E16').Value = '  -3.46%  '
$ws.Range('D17').Value = "'14.16"
$ws.Range('E17').Value = '  -3.79%  '
$ws.Range('D18').Value = '43.206.90'
$ws.Range('E18').Value = '  -5.97%  '
$ws.Range('D19').Value = "'13.24"
$ws.Range('E19').Value = '  +5.73%  '
$ws.Range('D20').Value = '0.0₃0983'
$ws.Range('E20').Value = '  -2.49%  '
$ws.Range('E21').Value = '  -0.56%  '
$ws.Range('D22').Value = "'72.14"
$ws.Range('E22').Value = '  -1.15%  '
$ws.Range('D23').Value = "'260.47"
$ws.Range('E23').Value = '  -10.48%  '
$ws.Range('D25').Value = "'29.71"
$ws.Range('E25').Value = '  +1.97%  '
$ws.Range('D26').Value = "'2.13"
$ws.Range('E26').Value = '  -4.84%  '
$ws.Range('E27').Value = '  -0.03%  '
$ws.Range('D28').Value = "'10.06"
$ws.Range('E28').Value = '  -5.62%  '
$ws.Range('D29').Value = "'37.71"
$ws.Range('E29').Value = '  -1.56%  '
$ws.Range('E30').Value = '  -5.59%  '
$ws.Range('E31').Value = '  -4.22%  '
$ws.Range('D32').Value = "'154.45"
$ws.Range('E32').Value = '  -2.38%  '
$ws.Range('D33').Value = "'2.19"
$ws.Range('E33').Value = '  -0.95%  '
$ws.Range('E34').Value = '  -1.66%  '
$ws.Range('E35').Value = '  -6.85%  '
$ws.Range('D36').Value = "'0.0800"
$ws.Range('E36').Value = '  -4.11%  '
$ws.Range('D37').Value = "'0.116"
$ws.Range('E37').Value = '  -3.90%  '
$ws.Range('E38').Value = '  -2.11%  '
$ws.Range('D39').Value = "'17.04"
$ws.Range('E39').Value = '  +9.99%  '
$ws.Range('D40').Value = "'23.46"
$ws.Range('E40').Value = '  +11.09%  '
$ws.Range('E41').Value = '  -1.41%  '
$ws.Range('D42').Value = "'0.0313"
$ws.Range('E42').Value = '  -3.91%  '
$ws.Range('E43').Value = '  -1.82%  '
$ws.Range('D44').Value = '2.085.51'
$ws.Range('E44').Value = '  -1.76%  '
$ws.Range('E45').Value = '  -0.12%  '
$ws.Range('D46').Value = "'86.06"
$ws.Range('E46').Value = '  -9.17%  '
$ws.Range('D47').Value = "'8.85"
$ws.Range('E47').Value = '  -4.19%  '
$ws.Range('E48').Value = '  +1.80%  '
$ws.Range('D49').Value = '2.804.39'
$ws.Range('E49').Value = '  -1.40%  '
$ws.Range('D50').Value = "'104.94"
$ws.Range('E50').Value = '  -3.15%  '
$ws.Range('E51').Value = '  -1.51%  '
